$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "21.01.2019"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Improved enemy AI."

[void]$ws.Range("C5").Select()
